$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New open/close/high/low/shares_outstanding/fixed_ticker values for rows 2-33
$data = @(
    @(13.79884986122726, 12.93764495849609, 14.43496666287938, 11.55776015903448, 708442998, "TME"),
    @(13.79884986122726, 12.93764495849609, 14.43496666287938, 11.55776015903448, 708442998, "TME"),
    @(13.79884986122726, 12.93764495849609, 14.43496666287938, 11.55776015903448, 708442998, "TME"),
    @(13.79884986122726, 12.93764495849609, 14.43496666287938, 11.55776015903448, 708442998, "TME"),
    @(13.79884986122726, 12.93764495849609, 14.43496666287938, 11.55776015903448, 708442998, "TME"),
    @(13.79884986122726, 12.93764495849609, 14.43496666287938, 11.55776015903448, 708442998, "TME"),
    @(12.50704165715845, 14.63069438934326, 15.8637823198816, 11.80242024922384, 708442998, "TME"),
    @(18.34953373644923, 16.83263969421387, 19.00522381662946, 16.42161006385296, 708442998, "TME"),
    @(15.0808705973974, 13.96521949768066, 16.38246416927066, 13.62269448053272, 708442998, "TME"),
    @(12.48746997446418, 13.54440307617188, 14.18051989233288, 12.29174179513046, 708442998, "TME"),
    @(11.77306092887055, 12.41896438598633, 13.90650041245165, 11.64583737924084, 708442998, "TME"),
    @(9.620048726547701, 11.16630268096924, 11.29841954668827, 9.531970816068348, 708442998, "TME"),
    @(13.33888826360324, 15.79527854919434, 17.58619302459886, 13.03550888545076, 708442998, "TME"),
    @(14.60133569040702, 14.56219005584717, 14.81663714713896, 13.55418926595188, 708442998, "TME"),
    @(18.88778566128825, 26.0318717956543, 28.22402916512052, 18.49632932149269, 708442998, "TME"),
    @(20.7080592088494, 17.04793930053711, 20.92335972112271, 16.46564813328663, 708442998, "TME"),
    @(15.16894859794533, 10.3442440032959, 15.28638550735211, 9.306884170216575, 708442998, "TME"),
    @(7.065793304166912, 7.692124366760254, 8.905640392213693, 6.684123148305763, 708442998, "TME"),
    @(6.772201220632248, 6.048006057739258, 7.496395916872471, 5.245520093435526, 708442998, "TME"),
    @(5.088937415662525, 4.159227848052979, 5.426569114065334, 3.728625379550335, 708442998, "TME"),
    @(4.863849406350651, 4.120081901550293, 5.274879025663927, 3.983072183996789, 708442998, "TME"),
    @(3.934139974783652, 3.532896757125854, 4.50175213798012, 3.072935321427852, 708442998, "TME"),
    @(8.269522923912934, 8.21080493927002, 9.091582175358058, 7.613833317784739, 708442998, "TME"),
    @(8.06400831763397, 7.251735687255859, 8.073794959760349, 6.742841963003743, 708442998, "TME"),
    @(7.418105192701361, 6.840706348419189, 7.77041639147123, 6.175230054434152, 708442998, "TME"),
    @(6.194802733472838, 7.104939937591553, 7.212590203322728, 5.959928917347229, 708442998, "TME"),
    @(8.636513280336795, 9.19923210144043, 9.37538791557237, 8.015075543420679, 708442998, "TME"),
    @(11.0195068818178, 12.28195476531982, 12.65383876097309, 10.70634087145921, 708442998, "TME"),
    @(13.99373191921562, 14.00360774993896, 15.3368138178252, 13.63821048947484, 708442998, "TME"),
    @(11.1199307861233, 11.8309736251831, 12.25562491597987, 10.13237034561809, 708442998, "TME"),
    @(14.26037246230983, 13.25306129455566, 14.3936933488954, 11.5643328899048, 708442998, "TME"),
    @(19.42000007629395, 20.98999977111816, 22.5, 18.82200050354004, 708442998, "TME")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Range("D$r").Value = $row[0]
    $ws.Range("E$r").Value = $row[1]
    $ws.Range("F$r").Value = $row[2]
    $ws.Range("G$r").Value = $row[3]
    $ws.Range("H$r").Value = $row[4]
    $ws.Range("I$r").Value = $row[5]
}
